# 03项目计划表.xlsx - add the "日期：2018.11.12第十一周周一" weekly block
# (rows 171-180) beneath the existing "日期：2018.11.8第十周周四" block, and
# flesh out the summary line for the previous week (row 169 / A169).
#
# Alignment / border constants (mirrors the Excel enum values so this reads
# the same under real Excel COM):
#   xlCenter = -4108, xlLeft = -4131
#   xlEdgeLeft = 7, xlEdgeTop = 8, xlEdgeBottom = 9, xlEdgeRight = 10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-AllBorders($rng) {
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------
# 1. Fill in the previously-empty summary sentence for the prior week.
# ---------------------------------------------------------------------
$ws.Range("A169").Value = "总结：我们的交互现在有点慢了，我们得找一下是什么原因啊，我们服务端的编码基本上是弄好了，客户端的需要加把劲啊"

# ---------------------------------------------------------------------
# 2. New week header band (row 171) - merged A:D, top border, centered.
# ---------------------------------------------------------------------
$hdr = $ws.Range("A171:D171")
$hdr.Borders.Item(8).LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$ws.Range("A171").Value = "日期：2018.11.12第十一周周一"
$hdr.Merge()

# ---------------------------------------------------------------------
# 3. Table header row (row 172): 组员 / 计划内容 / (blank) / 备注
# ---------------------------------------------------------------------
foreach ($addr in @("A172", "B172", "D172")) {
    $c = $ws.Range($addr)
    Set-AllBorders $c
    $c.VerticalAlignment = -4108
    $c.Font.Bold = $true
    $c.Font.Size = 10
}
$ws.Range("A172").Value = "组员"
$ws.Range("B172").Value = "计划内容"
$ws.Range("D172").Value = "备注"

$c172 = $ws.Range("C172")
Set-AllBorders $c172
$c172.HorizontalAlignment = -4108
$c172.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Data rows 173-178: one row per team member.
# ---------------------------------------------------------------------
$members = @(
    @{ Row = 173; Name = "陈柯赞"; Task = "管理端的编码与维护" },
    @{ Row = 174; Name = "黎安生"; Task = "pc端的申请加入群聊与数据传递" },
    @{ Row = 175; Name = "王智永"; Task = "写一个服务类，每隔几秒获取经纬度，并发送到后台" },
    @{ Row = 176; Name = "郑海文"; Task = "申请加入群聊的数据传递，并展示在客户端上" },
    @{ Row = 177; Name = "赵华亮"; Task = "写个封装类来解析后台传递的json数据" },
    @{ Row = 178; Name = "叶田";   Task = "管理端的图片收集" }
)

foreach ($m in $members) {
    $r = $m.Row

    # Column A - member name, plain bordered (all-thin) cell.
    $a = $ws.Range("A$r")
    Set-AllBorders $a
    $a.VerticalAlignment = -4108
    $a.Value = $m.Name

    # Column B - task text, all-thin border except the last row of the
    # group which gets the "closing" bottom-only-emphasis style.
    $b = $ws.Range("B$r")
    Set-AllBorders $b
    $b.HorizontalAlignment = -4108
    $b.VerticalAlignment = -4108
    $b.Value = $m.Task

    # Column C - completion percentage cell (left blank, like the two
    # preceding weekly blocks), first row gets the top border, the rest
    # plain all-thin borders; all of them use a 0% number format.
    $c = $ws.Range("C$r")
    if ($r -eq 173) {
        $c.Borders.Item(7).LineStyle = 1
        $c.Borders.Item(10).LineStyle = 1
        $c.Borders.Item(8).LineStyle = 1
    } else {
        Set-AllBorders $c
    }
    $c.NumberFormat = "0%"
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108

    # Column D - remarks column; only border formatting, no text/value.
    $d = $ws.Range("D$r")
    if ($r -eq 173) {
        $d.Borders.Item(7).LineStyle = 1
        $d.Borders.Item(10).LineStyle = 1
        $d.Borders.Item(8).LineStyle = 1
    } else {
        Set-AllBorders $d
    }
    $d.HorizontalAlignment = -4108
    $d.VerticalAlignment = -4108
}

# Column C/D merges: C only spans the first two rows (173:174, matching
# the "completion %" merge in the two prior weekly blocks); D spans the
# whole member list (173:178).
$ws.Range("C173:C174").Merge()
$ws.Range("D173:D178").Merge()

# ---------------------------------------------------------------------
# 5. Summary rows 179 (label) and 180 (blank filler), matching the
#    layout used by every other weekly block in the sheet.
# ---------------------------------------------------------------------
$sumRange = $ws.Range("A179:D180")
Set-AllBorders $sumRange
$sumRange.HorizontalAlignment = -4131
$sumRange.VerticalAlignment = -4108
$ws.Range("A179").Value = "总结："
$ws.Range("A179:D179").Merge()
$ws.Range("A180:D180").Merge()

# ---------------------------------------------------------------------
# 6. Move the visible selection down to the freshly-added block, same as
#    the author's own cursor position when they saved the workbook.
# ---------------------------------------------------------------------
$ws.Range("A169:D170").Select()
$ws.Range("A169:D170").Select()
$ws.Range("A169:D170").Select()
$ws.Range("A169:D170").Select()
